$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Rows 37 and 38 swap their entire record content (same field values moved
# from one row to the other). Only the cells that actually differ between
# the two rows are touched below.
# ---------------------------------------------------------------------------

# --- Row 37 becomes what Row 38 used to hold -------------------------------
$ws.Range("A37").Value = 131106643
$ws.Range("B37").Value = 78648
$ws.Range("E37").Value = 6437
$ws.Range("F37").Value = "Blanksvart spiklav"
$ws.Range("G37").Value = "Calicium denigratum"
$ws.Range("H37").Value = "(Vain.) Tibell"

# I37 needs to hold the text "1" (not the number 1), matching the source
# data's string-typed column. Force text formatting first so COM stores it
# as a string instead of auto-coercing it to a number.
$ws.Range("I37").NumberFormat = "@"
$ws.Range("I37").Value = "1"

$ws.Range("J37").Value = "cm²"
$ws.Range("Q37").Value = 601129
$ws.Range("R37").Value = 6959679
$ws.Range("X37").Value = "2025_0533"
$ws.Range("Z37").Value = "09:55"
$ws.Range("AB37").Value = "09:55"
$ws.Range("AC37").ClearContents()
$ws.Range("AX37").Value = "Samuel Koont"

# --- Row 38 becomes what Row 37 used to hold -------------------------------
$ws.Range("A38").Value = 131106646
$ws.Range("B38").Value = 79245
$ws.Range("E38").Value = 6425
$ws.Range("F38").Value = "Garnlav"
$ws.Range("G38").Value = "Alectoria sarmentosa"
$ws.Range("H38").Value = "(Ach.) Ach."
$ws.Range("I38").Value = ""
$ws.Range("J38").ClearContents()
$ws.Range("Q38").Value = 601173
$ws.Range("R38").Value = 6959739
$ws.Range("X38").Value = "2025_0530"
$ws.Range("Z38").Value = "09:36"
$ws.Range("AB38").Value = "09:36"
$ws.Range("AC38").Value = "tall"
$ws.Range("AX38").Value = "Måns Svensson"

# ---------------------------------------------------------------------------
# Small numeric corrections on rows 39 and 40
# ---------------------------------------------------------------------------
$ws.Range("B39").Value = 91824
$ws.Range("B40").Value = 98935
